# Generate Report for Handback
#
# The handback transform for file ff40acbf-6728-42b1-bf0f-62301f78ecef failed
# because the produced handback file name did not match the expected handoff
# file name. Update the status on the Overview sheet (and on each per-locale
# sheet's Status column) from "Ready for handoff" to "Handback transform
# failed", and record the error detail for each locale in column L
# ("Error Detail") of its row.

$wb = $excel.ActiveWorkbook

$newStatus = "Handback transform failed"

# --- Overview sheet: row for ff40acbf-... (row 3) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $newStatus   # zh-cn column
$overview.Range("C3").Value = $newStatus   # de-de column

# --- zh-cn sheet: row for ff40acbf-... (row 3) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("L3").Value = "Handback file name: 1gncqslm.gj1 is different with handoff file name: ff40acbf-6728-42b1-bf0f-62301f78ecef.3a26b5ca79b94dc00df8e783f2127b01f3823c06.zh-cn."

# --- de-de sheet: row for ff40acbf-... (row 3) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $newStatus
$dede.Range("L3").Value = "Handback file name: 1gncqslm.gj1 is different with handoff file name: ff40acbf-6728-42b1-bf0f-62301f78ecef.3a26b5ca79b94dc00df8e783f2127b01f3823c06.de-de."
